$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update June total for 2021 (H7)
$ws.Range("H7").Value = 131

# Update August row label to reflect new "through" date
$ws.Range("A9").Value = "August (through 08-20)"

# Update August counts per year
$ws.Range("B9").Value = 22
$ws.Range("C9").Value = 46
$ws.Range("D9").Value = 53
$ws.Range("G9").Value = 124
$ws.Range("H9").Value = 102

# Update Total row counts per year
$ws.Range("B10").Value = 184
$ws.Range("C10").Value = 348
$ws.Range("D10").Value = 518
$ws.Range("G10").Value = 745
$ws.Range("H10").Value = 1017
